$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "31.415.96"
$ws.Cells.Item(2, 5).Value = "  +5.18%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "1.706.17"
$ws.Cells.Item(3, 5).Value = "  +4.13%  "
# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.20%  "
# Row 5
$ws.Cells.Item(5, 4).Value = "221.69"
$ws.Cells.Item(5, 5).Value = "  +2.94%  "
# Row 6
$ws.Cells.Item(6, 5).Value = "  +2.86%  "
# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.18%  "
# Row 8
$ws.Cells.Item(8, 4).Value = "29.87"
$ws.Cells.Item(8, 5).Value = "  +4.04%  "
# Row 9
$ws.Cells.Item(9, 2).Value = "OKB"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(9, 4).Value = "45.38"
$ws.Cells.Item(9, 5).Value = "  +3.44%  "
# Row 10
$ws.Cells.Item(10, 2).Value = "Cardano"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(10, 4).Value = "0.269"
$ws.Cells.Item(10, 5).Value = "  +3.28%  "
# Row 11
$ws.Cells.Item(11, 2).Value = "Dogecoin"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(11, 4).Value = "0.0643"
$ws.Cells.Item(11, 5).Value = "  +5.59%  "
# Row 12
$ws.Cells.Item(12, 2).Value = "TRON"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(12, 4).Value = "0.0911"
$ws.Cells.Item(12, 5).Value = "  +1.24%  "
# Row 13
$ws.Cells.Item(13, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(13, 4).Value = "1.951.60"
$ws.Cells.Item(13, 5).Value = "  +4.27%  "
# Row 14
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "1.712.60"
$ws.Cells.Item(14, 5).Value = "  +4.38%  "
# Row 15
$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).Value = "10.28"
$ws.Cells.Item(15, 5).Value = "  +9.23%  "
# Row 16
$ws.Cells.Item(16, 2).Value = "Polygon"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(16, 4).Value = "0.611"
$ws.Cells.Item(16, 5).Value = "  +3.13%  "
# Row 17
$ws.Cells.Item(17, 2).Value = "Polkadot"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(17, 4).Value = "4.18"
$ws.Cells.Item(17, 5).Value = "  +8.36%  "
# Row 18
$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(18, 4).Value = "31.369.44"
$ws.Cells.Item(18, 5).Value = "  +4.99%  "
# Row 19
$ws.Cells.Item(19, 2).Value = "Litecoin"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(19, 4).Value = "67.08"
$ws.Cells.Item(19, 5).Value = "  +3.78%  "
# Row 20
$ws.Cells.Item(20, 2).Value = "BitcoinCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(20, 4).Value = "250.77"
$ws.Cells.Item(20, 5).Value = "  +4.32%  "
# Row 21
$ws.Cells.Item(21, 2).Value = "ShibaInu"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(21, 4).Value = "0.0₃0725"
$ws.Cells.Item(21, 5).Value = "  +3.15%  "
# Row 22
$ws.Cells.Item(22, 2).Value = "Dai"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(22, 4).Value = "0.998"
$ws.Cells.Item(22, 5).Value = "  -0.20%  "
# Row 23
$ws.Cells.Item(23, 2).Value = "Avalanche"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(23, 4).Value = "10.20"
$ws.Cells.Item(23, 5).Value = "  +2.95%  "
# Row 24
$ws.Cells.Item(24, 2).Value = "Uniswap"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(24, 4).Value = "4.28"
$ws.Cells.Item(24, 5).Value = "  +3.40%  "
# Row 25
$ws.Cells.Item(25, 2).Value = "Toncoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(25, 4).Value = "2.16"
$ws.Cells.Item(25, 5).Value = "  -0.38%  "
# Row 26
$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).Value = "159.23"
$ws.Cells.Item(26, 5).Value = "  +1.06%  "
# Row 27
$ws.Cells.Item(27, 2).Value = "EthereumClassic"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).Value = "16.03"
# Row 28
$ws.Cells.Item(28, 2).Value = "Stellar"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(28, 4).Value = "0.113"
$ws.Cells.Item(28, 5).Value = "  +2.99%  "
# Row 29
$ws.Cells.Item(29, 2).Value = "Cosmos"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(29, 4).Value = "6.80"
$ws.Cells.Item(29, 5).Value = "  +2.67%  "
# Row 30
$ws.Cells.Item(30, 2).Value = "BinanceUSD"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(30, 4).Value = "0.999"
$ws.Cells.Item(30, 5).Value = "  -0.16%  "
# Row 31
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).Value = "3.79"
$ws.Cells.Item(31, 5).Value = "  +11.64%  "
# Row 32
$ws.Cells.Item(32, 2).Value = "Hedera"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(32, 4).Value = "0.0504"
$ws.Cells.Item(32, 5).Value = "  +2.37%  "
# Row 33
$ws.Cells.Item(33, 2).Value = "PancakeSwap"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(33, 4).Value = "1.15"
$ws.Cells.Item(33, 5).Value = "  +3.64%  "
# Row 34
$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).Value = "3.42"
$ws.Cells.Item(34, 5).Value = "  +7.01%  "
# Row 35
$ws.Cells.Item(35, 2).Value = "Maker"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(35, 4).Value = "1.511.81"
$ws.Cells.Item(35, 5).Value = "  +6.07%  "
# Row 36
$ws.Cells.Item(36, 2).Value = "LidoDAOToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(36, 4).Value = "1.74"
$ws.Cells.Item(36, 5).Value = "  +2.41%  "
# Row 37
$ws.Cells.Item(37, 2).Value = "TrustWalletToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(37, 4).Value = "1.04"
$ws.Cells.Item(37, 5).Value = "  +1.73%  "
# Row 38
$ws.Cells.Item(38, 2).Value = "Aave"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(38, 4).Value = "83.44"
$ws.Cells.Item(38, 5).Value = "  +9.33%  "
# Row 39
$ws.Cells.Item(39, 2).Value = "ImmutableX"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(39, 4).Value = "0.610"
$ws.Cells.Item(39, 5).Value = "  +8.91%  "
# Row 40
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(40, 4).Value = "0.0181"
$ws.Cells.Item(40, 5).Value = "  +4.64%  "
# Row 41
$ws.Cells.Item(41, 2).Value = "MXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(41, 4).Value = "2.72"
$ws.Cells.Item(41, 5).Value = "  -0.46%  "
# Row 42
$ws.Cells.Item(42, 2).Value = "HuobiToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(42, 4).Value = "2.30"
$ws.Cells.Item(42, 5).Value = "  +0.45%  "
# Row 43
$ws.Cells.Item(43, 2).Value = "RenderToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(43, 4).Value = "2.04"
$ws.Cells.Item(43, 5).Value = "  +3.71%  "
# Row 44
$ws.Cells.Item(44, 2).Value = "ARBITRUM"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(44, 4).Value = "0.854"
$ws.Cells.Item(44, 5).Value = "  +2.24%  "
# Row 45
$ws.Cells.Item(45, 2).Value = "Kaspa"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(45, 4).Value = "0.0502"
$ws.Cells.Item(45, 5).Value = "  +0.40%  "
# Row 46
$ws.Cells.Item(46, 2).Value = "WEMIXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(46, 4).Value = "1.04"
$ws.Cells.Item(46, 5).Value = "  +3.20%  "
# Row 47
$ws.Cells.Item(47, 2).Value = "PaxDollar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(47, 4).Value = "0.999"
$ws.Cells.Item(47, 5).Value = "  -0.11%  "
# Row 48
$ws.Cells.Item(48, 2).Value = "BitcoinSV"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(48, 4).Value = "52.07"
$ws.Cells.Item(48, 5).Value = "  +6.32%  "
# Row 49
$ws.Cells.Item(49, 2).Value = "FraxShare"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(49, 4).Value = "5.59"
$ws.Cells.Item(49, 5).Value = "  +3.67%  "
# Row 50
$ws.Cells.Item(50, 2).Value = "RocketPoolETH"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(50, 4).Value = "1.843.83"
$ws.Cells.Item(50, 5).Value = "  +3.65%  "
# Row 51
$ws.Cells.Item(51, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51, 4).Value = "0.0₆0119"
$ws.Cells.Item(51, 5).Value = "  +8.97%  "
